$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 843, shifting existing rows 843:941 down to 846:944.
$ws.Rows("843:845").Insert()

# Populate the 3 newly inserted rows (843, 844, 845) with new weekly records.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are constant across this dataset.

# Row 843 - Especial
$ws.Range("A843").Value = 8
$ws.Range("B843").Value = "Terminal La Palmera de La Serena"
$ws.Range("C843").Value = "Coquimbo"
$ws.Range("D843").Value = 44918
$ws.Range("E843").Value = 4
$ws.Range("F843").Value = "Fruta"
$ws.Range("G843").Value = 100101
$ws.Range("H843").Value = "Berries"
$ws.Range("I843").Value = 100112025
$ws.Range("J843").Value = "Frutilla"
$ws.Range("K843").Value = "Sin especificar"
$ws.Range("L843").Value = "Especial"
$ws.Range("M843").Value = 400
$ws.Range("N843").Value = 12000
$ws.Range("O843").Value = 13000
$ws.Range("P843").Value = 12500
$ws.Range("Q843").Value = "`$/bandeja 7 kilos"
$ws.Range("R843").Value = "Provincia de Melipilla"
$ws.Range("S843").Value = 1786
$ws.Range("T843").Value = 7

# Row 844 - Primera
$ws.Range("A844").Value = 8
$ws.Range("B844").Value = "Terminal La Palmera de La Serena"
$ws.Range("C844").Value = "Coquimbo"
$ws.Range("D844").Value = 44918
$ws.Range("E844").Value = 4
$ws.Range("F844").Value = "Fruta"
$ws.Range("G844").Value = 100101
$ws.Range("H844").Value = "Berries"
$ws.Range("I844").Value = 100112025
$ws.Range("J844").Value = "Frutilla"
$ws.Range("K844").Value = "Sin especificar"
$ws.Range("L844").Value = "Primera"
$ws.Range("M844").Value = 400
$ws.Range("N844").Value = 10000
$ws.Range("O844").Value = 11000
$ws.Range("P844").Value = 10500
$ws.Range("Q844").Value = "`$/bandeja 7 kilos"
$ws.Range("R844").Value = "Provincia de Melipilla"
$ws.Range("S844").Value = 1500
$ws.Range("T844").Value = 7

# Row 845 - Segunda
$ws.Range("A845").Value = 8
$ws.Range("B845").Value = "Terminal La Palmera de La Serena"
$ws.Range("C845").Value = "Coquimbo"
$ws.Range("D845").Value = 44918
$ws.Range("E845").Value = 4
$ws.Range("F845").Value = "Fruta"
$ws.Range("G845").Value = 100101
$ws.Range("H845").Value = "Berries"
$ws.Range("I845").Value = 100112025
$ws.Range("J845").Value = "Frutilla"
$ws.Range("K845").Value = "Sin especificar"
$ws.Range("L845").Value = "Segunda"
$ws.Range("M845").Value = 500
$ws.Range("N845").Value = 8000
$ws.Range("O845").Value = 9000
$ws.Range("P845").Value = 8500
$ws.Range("Q845").Value = "`$/bandeja 7 kilos"
$ws.Range("R845").Value = "Provincia de Melipilla"
$ws.Range("S845").Value = 1214
$ws.Range("T845").Value = 7

# Ensure the date cells carry the same date-formatted style as the rest of column D.
$ws.Range("D846").Copy()
$ws.Range("D843:D845").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D843").Value = 44918
$ws.Range("D844").Value = 44918
$ws.Range("D845").Value = 44918
